$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (post-edit) player/position/team table, rows 2-19 (row 1 is the header).
$data = @(
    @("Donovan Mitchell",    "PG,SG",    "Cleveland Cavaliers"),
    @("Ty Jerome",           "PG,SG",    "Cleveland Cavaliers"),
    @("Malik Beasley",       "SG,SF",    "Detroit Pistons"),
    @("Dyson Daniels",       "PG,SG,SF", "Atlanta Hawks"),
    @("Jaden McDaniels",     "SF,PF",    "Minnesota Timberwolves"),
    @("Michael Porter Jr.",  "SF,PF",    "Denver Nuggets"),
    @("Kristaps Porzingis",  "PF,C",     "Boston Celtics"),
    @("Josh Hart",           "SG,SF,PF", "New York Knicks"),
    @("Kel'el Ware",         "PF,C",     "Miami Heat"),
    @("Alperen Sengün",      "C",        "Houston Rockets"),
    @("Domantas Sabonis",    "C",        "Sacramento Kings"),
    @("Victor Wembanyama",   "C",        "San Antonio Spurs"),
    @("Kelly Oubre Jr.",     "SG,SF",    "Philadelphia 76ers"),
    @("Andrew Wiggins",      "SF,PF",    "Golden State Warriors"),
    @("Guerschon Yabusele",  "PF,C",     "Philadelphia 76ers"),
    @("Cam Thomas",          "SG,SF",    "Brooklyn Nets"),
    @("Donte DiVincenzo",    "PG,SG,SF", "Minnesota Timberwolves"),
    @("De'Andre Hunter",     "SF,PF",    "Atlanta Hawks")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
